# Identity and report fix
# Append a new data row (row 72) to the comments report sheet:
#   A72 = 70 (sequential number)
#   B72 = 76 (id)
#   C72 = "test" (comment text - new shared string)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A72").Value = 70
$ws.Range("B72").Value = 76
$ws.Range("C72").Value = "test"
